# TestCase_Template.xlsx - add "Type" and "Value" columns to the test-case
# table, drop the two unused blank sheets and rename the remaining sheet.
# ("test data's document create function.")

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- drop the unused Sheet2 / Sheet3, keep only the data sheet ----------
$wb.Worksheets("Sheet2").Delete()
$wb.Worksheets("Sheet3").Delete()
$wb.Worksheets("Sheet1").Name = "TestCase"

$ws = $wb.Worksheets("TestCase")

# --- insert the two new columns ------------------------------------------
# Original layout: A=No. B=Element C=Action D=Screenshot(Yes/No)
# New layout:       A=No. B=Type C=Element D=Action E=Value F=Screenshot(Yes/No)
$ws.Columns("B").Insert()
$ws.Columns("E").Insert()

# --- header row -----------------------------------------------------------
$ws.Range("B4").Value = "Type"
$ws.Range("E4").Value = "Value"

# --- column widths ---------------------------------------------------------
$ws.Columns("E").ColumnWidth = 13.6640625

# --- fix up the conditional formatting ranges (the column insert does not
#     move these automatically, unlike data validation) -------------------
$cf1 = $ws.Range("D5:D10").FormatConditions
foreach ($fc in $cf1) {
    $fc.ModifyAppliesToRange($ws.Range("F5:F10"))
}
$cf2 = $ws.Range("D11:D33").FormatConditions
foreach ($fc in $cf2) {
    $fc.ModifyAppliesToRange($ws.Range("F11:F33"))
}

# --- restore the active selection ----------------------------------------
$ws.Range("H15").Select()
